# Replace orch_session_issue_id UUID values (column L) with freshly
# generated UUIDs, as produced by a new diagnostics run (observation #79
# interpretation fix).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$map = @{
    "18f35809-12c4-4686-bbbd-d3340cce7601" = "acd21a00-9511-48f4-80c3-cead3641eab5"
    "8f722592-e71d-4b05-8e7c-3b30bef89e70" = "a5176fc8-1b39-46ff-8a9e-3b2637ce576d"
    "34df525a-32e5-493c-a1f4-d3408ee4c4d7" = "ecfc54ed-a1b2-4a03-a474-570f0327e13f"
    "9469f810-0689-4da8-a980-621368c15427" = "cc483033-701a-464d-870a-121c406d995d"
    "d472a582-0b4b-4120-aed9-301962c39103" = "ee9cf181-0c27-491e-b428-4bbdb13abf68"
    "99b97686-23ad-4866-8dff-640d5ed58c09" = "e5f726d2-a08a-4017-82fb-6f8adc775ad7"
    "a5a1d6bc-2ea0-4cc6-92e3-089e1bd9f240" = "649bff15-140d-4fc5-9a4d-12374ec0536d"
    "7e822635-51ee-4273-8f34-a21a0784dd9b" = "d924d541-1f77-4a02-a05e-5b45713daef4"
    "dcdf57ff-b96c-4048-9640-bbf1db7f854d" = "0014ef4d-b4f9-4bd5-bf4a-58eb318ac552"
    "24e6e0e4-6c11-4ee1-9563-7c926f8e8790" = "b0292c27-661a-41f7-b50c-98a14d950b10"
    "7e505f38-0a8a-4f62-ac05-885d80f954fb" = "1fbbbdd0-a252-49b1-81e2-4cbdba850ba3"
    "4444bdac-5b7c-477f-9753-084dbdcd4494" = "9c30ee6e-4f96-4361-85e8-c091c11e333b"
    "2eb16b50-6b3b-48ce-9cca-c9faed6f1cae" = "b51e797e-6a03-4643-ad4f-bcee158c9f1b"
    "0f976ad7-b6ef-4794-a7a9-898954b41f6f" = "6e9f267a-20c6-42a8-9a76-cd9c05265110"
    "63caee03-619a-45e0-9b8b-adc8582bcd94" = "debeaf96-c691-4074-9c85-086be4208b84"
    "62f5530c-c0db-4c70-bff3-9ec6eda6e377" = "f54a0b04-8459-472f-91cc-662eea6c7680"
    "7f3ae49e-8eca-44cf-ad21-5db6632fc78c" = "5d29a9e6-a8a5-447c-a52f-036c174e52ab"
    "0a20d09e-6a2e-4344-a426-810035742d85" = "567dd62c-9034-4f3b-b4b3-a729c3c8e08e"
    "19c1da9c-6315-4dc3-8a3a-3233bf564ac7" = "aabf14f0-257f-43c9-a4fc-278a346ffa58"
    "0db7fcac-77b7-4a7f-935a-345e23e9c06b" = "aab2cf71-dca1-452a-bf14-b56062d9d4f1"
    "8c076cf1-86de-4ced-9141-c0c840179a1c" = "738e042b-2620-4c82-a642-a5beb55a6a02"
    "34dc9014-4813-4f8d-895d-4cc2c4b67ecd" = "1c8289ed-7608-4e6b-b7b6-082a8ff11ecf"
    "d3548573-63d9-4acb-9f4d-62480f2f062f" = "1aec4b08-791c-40ae-aa4e-099ba2df7423"
    "40ccd34a-ee04-42d0-88bf-9924dcce854f" = "8fa111e5-8af2-4396-9970-a746ee376c85"
    "6e724ffc-e219-4a40-8274-99d395a9535d" = "f0bf1e50-97bb-4706-9b9f-56b0cfb20bb4"
    "64045e52-a909-480a-8f12-c0061b176955" = "47ca3904-83ea-435c-9b67-5cd16c013598"
    "d76bab4a-fdd5-4fbf-b79a-933e0c7671a7" = "c79778a9-cd6c-4572-ab3a-2bd0c469d359"
    "68c2108d-5bb5-4e34-b61f-b79c5d70fa26" = "1903820c-3c30-4de0-b42f-6e9d27440810"
    "d48bed29-3aa3-470d-be47-9ad4c18faebd" = "8518be33-0bb4-499b-aba0-c9dabb89a160"
    "43a1abc9-008f-4814-81a1-71913afd2058" = "4ef19dd6-1c73-40e1-af0d-30f609cbe595"
    "5e06ffda-e146-47b2-9952-13162bc2bc03" = "6fd52d2e-5a6c-4ee5-8522-393aa7877a87"
    "94f35b70-8d70-4fd5-b908-5633abda107d" = "97a205b2-98f5-4fd7-8aef-c17415d402ec"
    "f95afc54-2829-4f35-8e00-c19f919bdec2" = "693e6a04-3be5-494c-a9b1-67a69e30f51a"
    "709eec19-a26d-45ac-9753-d3a33b8e797d" = "ac297eb7-8ac6-4052-acb2-4709ef7e0a31"
    "f7789803-62b8-4452-a9e1-d3572387d7d1" = "c3879b21-8ef5-4ffe-97ca-c00a6c45d849"
    "f6fed19c-f5a7-4374-b3c6-258c13e874ce" = "7bd429b0-eae0-42fe-8eec-7d21ef655dfd"
    "236a2fda-352f-43ad-b215-7eb27b90aadf" = "5b5d7707-1f12-4d58-a7ca-56e85111ca23"
    "2f69b10c-cc49-4409-80bc-082acc5c4f94" = "06b68a85-7976-4fa2-ba6f-a121eb363fca"
    "5cb54e2e-d74a-40f9-b3d4-1c89d634192a" = "3dbd6b1d-ff17-40e6-8afc-81d4ef1331ea"
    "1f7c0114-dadb-463e-8e44-1fbe90dd1b19" = "32afeb54-8ad6-4827-b970-0046d04a7cc8"
    "9d2dbdfa-baa8-4b8f-b757-b2b7b6e859aa" = "dd6dff5e-3389-4b43-8f65-7cb773adfb40"
    "1cf085f6-f7f2-4a9a-9103-82afdbc40947" = "26851362-3e2d-4643-8dfe-765ab4e73f78"
    "3c1521f8-b6e4-427c-a394-db4da2243e49" = "c14f5105-f52e-45d5-8ddf-f51df4ff658d"
    "282f9cca-a8d2-4563-93cf-0cf65e5862a4" = "a48990c9-04c9-438c-a543-efe791b77dd7"
    "ac11819e-1046-4c6b-9cee-811a34370e23" = "6949647c-8eff-42de-9034-764d207e1b29"
    "d87c1376-2c74-4ecb-8be1-d299e418fb9f" = "c82edfd8-1811-4947-ba2a-9447b2d7b77e"
    "11939701-a721-4067-b850-70c0dd57d1ec" = "864f006f-3daf-407f-b3be-d7b4d9e795d2"
    "a0d31bb3-d1d2-4ba0-a1ea-e4bd750b3203" = "e963b326-b443-4727-9cc1-aca1a9ef2f97"
    "dc1b16bd-486e-437e-af94-42d64fc92142" = "09811053-03e8-4351-b9ed-927889907f6c"
    "fea8e7fe-4105-4b7d-baba-0edae95646fc" = "5c4657dc-8cff-4f78-83f8-8356e608aada"
    "ddf2813a-fdb5-4e49-a5f3-f7f2072248ff" = "a73e0c04-7f06-400b-91f8-4a9d4553d99c"
    "7119a7e7-526f-4987-bb34-7970ea26bf71" = "9b213f1a-2551-4415-9b20-b98782006e5e"
    "8a1b8e04-b5d7-483c-ab60-c1f3410bc1d5" = "7a31426b-9330-4387-b9f0-a56a0603310a"
    "6ef5e26d-765d-47a9-9999-21e06c523cac" = "edb2c003-6fd8-4b62-b87f-36be0eca2ab2"
    "f0ea9208-cd26-46f9-834e-45a440999929" = "1611130b-33a0-4b1c-8e4e-694348e4a344"
    "c5d16bf4-bb77-43cd-812a-6d03dfa11f7f" = "a5505694-0081-4d4d-9da9-dc6eb0e4589c"
    "b4c3394d-c042-442a-b2a7-11431727dc25" = "c443079e-1771-4ccf-aeee-ba5f2ea821b7"
    "929c2fcb-5731-44f2-9eb8-038715fd63fc" = "6e63d31e-0fdb-41a7-bbf5-f2a9a005db6e"
    "f0f97103-6f38-454e-bc3e-40ce70a0b180" = "0710b8b5-7b25-445e-a118-ab9bbbf436ee"
    "9d15cc56-53a0-4f08-8f71-84d9ace82416" = "ef391bf4-ad8e-488c-b4a6-9e92f7b612d1"
    "6ef391fa-a23f-455c-8dab-8064d32bf49a" = "801fa8dd-c990-4b62-9672-35b863b0b54d"
    "69a287c0-495c-436c-9baa-ecbb3c0c6ed9" = "9b2f456d-f1ba-4638-a38c-d4d0add84488"
    "c75c1085-b113-44f2-899f-00dc17008432" = "79d3c887-6111-4ee2-bb15-356e02d5d5fa"
    "fb1122ca-f6d0-484d-8b96-749b2c8a2ca0" = "6e176bde-3792-42f5-9d5c-45bf3c40ed0a"
    "4271086d-1f1f-4de4-b2e3-163b0b5e7394" = "558d1955-2fb8-4aea-ab86-a6dfdedf5dde"
    "d64ad7a4-5762-4705-bf1e-b6231cf8a93d" = "35509e26-087f-4b6c-9aec-e6c0fb2599d9"
    "8853e600-59db-4081-8fd0-0cffca526a10" = "85699b23-56ed-424f-ab58-96d44ba3547d"
    "5c699d09-afbb-469b-b3d1-642593c37da1" = "1d2454c1-cb00-4203-aa0e-752a112c9d71"
    "cf02082d-2405-4ea6-970a-644bb7d43823" = "4566e1c7-b260-406e-9dc7-335d9b51e84a"
    "5394fdb4-8b40-4961-9202-90aae64dfc85" = "936c6634-cf5c-4a84-a515-9efc0d7558a7"
    "2a595655-da66-406d-935e-c15842bc23a2" = "5f7a94ff-8d8e-4414-8974-ad3f4f3bfe33"
    "87656cb0-b574-482f-ac26-30ecbbd14b53" = "0daea509-02a1-44ea-8f16-a6532f5aa2db"
    "89ba1e7a-413d-4c86-95d8-009fa86752ae" = "12b91a47-a951-44c6-9f5a-e4b64c1568bc"
    "7e47d60d-aa93-4e8f-b5b8-0d9aa55b98be" = "46c1984b-4591-4a2a-9f69-c6324793c3dd"
    "60fad860-feb1-44dd-9d68-204b54b96fdf" = "1bf4ddd3-062f-4ade-95f7-3577f5129124"
    "4923dc54-100a-4bd7-b9fc-5eb7867a65c3" = "6528c3f4-4a64-42ba-89e5-14342d8d794b"
    "a6e60b9f-f570-4a9e-a956-4445459fc626" = "d60704d3-536a-41e8-bb32-15ea1431e18a"
    "b895de91-8a9a-4acc-8aa3-2fe25f522190" = "0eda4177-5608-49d4-9db9-5eb72a7b91fc"
    "f4ece211-1a3e-4dfc-ad7c-2b3d51d9c8e8" = "d0c95e55-c098-4dbc-a9c0-8779ae476660"
    "0915beef-3bbc-441c-879f-a289a92af788" = "6e1f2434-21a5-4b11-aeb1-af482a57392f"
    "488135a0-bb4d-4ed6-83d5-f70b845fbee6" = "8037f2ef-e125-4bca-ac0e-be3ff1dc1474"
    "fed6f6fb-6623-447b-9097-33574a0ac4d1" = "59b9e6d9-e4e6-4636-8cc4-05d3854ffb21"
    "c79ea32f-dc13-4c16-a1e9-ad997ca2ec2b" = "2213cbde-f7cc-4943-bb6c-ba174cb0c4b5"
    "5dc2cd86-31e3-4cf3-b572-d50a6e61d877" = "56c2c6fc-4885-4142-afbd-d2f6aa719c90"
    "3006a651-8a3b-4361-bbe4-bfd120b6d1a3" = "12bee8c6-4dd7-45c5-8e55-eeb616020bf4"
    "72279112-c5c0-478f-b522-de27343bc6a8" = "3eb261de-9001-42e0-afb9-fd861cf16e75"
    "3a2231ab-3981-40f4-b4c5-9b535203c6aa" = "1f761cd3-d176-4b7c-a102-bdcae8cb9f44"
    "092db5f5-d166-4475-8bb0-7a98b1d4115d" = "378b1777-5a10-4d4a-a442-059ce3d7c082"
    "73668bd9-3aad-4721-b4cd-50991f01b913" = "ef76d150-9999-4877-a8b9-a142ae3f98a9"
    "4f5be7f8-dfde-40f2-9a2a-9276690cf0b8" = "66f68cd5-c981-4ba0-95f4-dbc2a55972d5"
    "0ff813c8-6ca6-40d3-90e1-8a83aba8fe4c" = "1dbba41f-4527-4a81-afb7-b92f6bc0740a"
    "566d1647-3598-4b8b-ad23-fa72324c716e" = "37d1d80d-e7d0-4b6b-a1eb-d85aed2f6a4b"
    "8a9ccda6-55ee-4185-a15e-4094c0cccd2a" = "3bb2e3e3-6f3f-4c45-96b8-552880b2c669"
    "ce68942c-1f9f-4753-8869-c8e4d2e6be80" = "4cd932df-49f6-42e5-9a66-5ccb03d5b1db"
    "92560696-f8a5-4120-9b32-a1a03de932d1" = "fc5f02b3-048a-49eb-8a52-2c03c256d4b6"
    "634118ee-3d99-47d1-9721-7d60d8ab06ae" = "2e256766-2bed-40ae-9feb-599a5d131a5c"
    "be5cd47c-eefc-4bd4-878c-030dadf60962" = "d8c10733-24f8-4d3d-aa76-1eb89b665be8"
    "34a504bb-63b1-483c-a99a-378da6092009" = "3e0aae8a-86c8-4b0b-925c-370dd8a8dedd"
    "11f01971-da61-48db-a254-621419a49660" = "5c48bb3c-10d0-4da2-845b-1f76bfc5ff91"
    "e00b7b88-7675-478f-8490-655a8c8cbd33" = "d456a5fe-9b10-4652-b169-5f0feacfef43"
    "7a70f0fc-6794-4521-b20e-e092d0f4b956" = "c12a0fb5-56f6-4bfb-8988-c40c2b976b47"
    "8c6da3ba-e511-4559-ad92-54518500e979" = "2d3bf0d1-5f03-4b83-92ec-d8f13d74ff78"
    "26c90ffb-d623-406f-8808-67dc499faed9" = "99d7158e-91a4-4ed8-825e-4a351a2d7fa9"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colIndex = 12  # column L = orch_session_issue_id

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $val = $cell.Value2
    if ($val -ne $null -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
